# Applies the market-price / profit figure updates described in the commit diff
# (Sheets/Alexander_Profits.xlsx) across the ALC, ARM, BSM, CRP, GSM, LTW and WVR tabs.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(124, 8).Value = 33822.223  # H124
$ws.Cells.Item(124, 10).Value = 33822.223  # J124
$ws.Cells.Item(124, 12).Value = 33822.223  # L124
$ws.Cells.Item(124, 14).Value = -43642.223  # N124
$ws.Cells.Item(126, 8).Value = 40595  # H126
$ws.Cells.Item(126, 10).Value = 40595  # J126
$ws.Cells.Item(126, 12).Value = 40595  # L126
$ws.Cells.Item(126, 14).Value = -50475  # N126
$ws.Cells.Item(128, 8).Value = 25442.857  # H128
$ws.Cells.Item(128, 10).Value = 25442.857  # J128
$ws.Cells.Item(128, 12).Value = 25442.857  # L128
$ws.Cells.Item(128, 14).Value = -35402.857  # N128
$ws.Cells.Item(137, 8).Value = 3948782.2  # H137
$ws.Cells.Item(137, 9).Value = 2084667.8  # I137
$ws.Cells.Item(137, 10).Value = 7144407.5  # J137
$ws.Cells.Item(137, 11).Value = 6254003.4  # K137
$ws.Cells.Item(137, 12).Value = 21433222.5  # L137
$ws.Cells.Item(137, 13).Value = -6251453.4  # M137
$ws.Cells.Item(137, 14).Value = -21438322.5  # N137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1043.1714  # H74
$ws.Cells.Item(74, 9).Value = 1062.8518  # I74
$ws.Cells.Item(74, 10).Value = 976.75  # J74
$ws.Cells.Item(74, 11).Value = 1062.8518  # K74
$ws.Cells.Item(74, 12).Value = 976.75  # L74
$ws.Cells.Item(74, 13).Value = -188.8517999999999  # M74
$ws.Cells.Item(74, 14).Value = -2724.75  # N74
$ws.Cells.Item(77, 8).Value = 1043.1714  # H77
$ws.Cells.Item(77, 9).Value = 1062.8518  # I77
$ws.Cells.Item(77, 10).Value = 976.75  # J77
$ws.Cells.Item(77, 11).Value = 5314.259  # K77
$ws.Cells.Item(77, 12).Value = 4883.75  # L77
$ws.Cells.Item(77, 13).Value = -946.259  # M77
$ws.Cells.Item(77, 14).Value = -13619.75  # N77
$ws.Cells.Item(80, 8).Value = 20495.5  # H80
$ws.Cells.Item(80, 10).Value = 20495.5  # J80
$ws.Cells.Item(80, 12).Value = 20495.5  # L80
$ws.Cells.Item(80, 14).Value = -22491.5  # N80
$ws.Cells.Item(83, 8).Value = 20495.5  # H83
$ws.Cells.Item(83, 10).Value = 20495.5  # J83
$ws.Cells.Item(83, 12).Value = 61486.5  # L83
$ws.Cells.Item(83, 14).Value = -71470.5  # N83
$ws.Cells.Item(109, 8).Value = 18320  # H109
$ws.Cells.Item(109, 10).Value = 18320  # J109
$ws.Cells.Item(109, 12).Value = 18320  # L109
$ws.Cells.Item(109, 14).Value = -21094  # N109
$ws.Cells.Item(127, 8).Value = 34863.57  # H127
$ws.Cells.Item(127, 10).Value = 34863.57  # J127
$ws.Cells.Item(127, 12).Value = 34863.57  # L127
$ws.Cells.Item(127, 14).Value = -44783.57  # N127
$ws.Cells.Item(132, 8).Value = 133237  # H132
$ws.Cells.Item(132, 9).Value = 162744.97  # I132
$ws.Cells.Item(132, 11).Value = 488234.91  # K132
$ws.Cells.Item(132, 13).Value = -485704.91  # M132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 3495.5  # H105
$ws.Cells.Item(105, 9).Value = 2123  # I105
$ws.Cells.Item(105, 11).Value = 2123  # K105
$ws.Cells.Item(105, 13).Value = -376  # M105
$ws.Cells.Item(122, 8).Value = 48500  # H122
$ws.Cells.Item(122, 10).Value = 48500  # J122
$ws.Cells.Item(122, 12).Value = 48500  # L122
$ws.Cells.Item(122, 14).Value = -58300  # N122
$ws.Cells.Item(132, 8).Value = 48500  # H132
$ws.Cells.Item(132, 10).Value = 48500  # J132
$ws.Cells.Item(132, 12).Value = 48500  # L132
$ws.Cells.Item(132, 14).Value = -58620  # N132

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 49099.6  # H20
$ws.Cells.Item(20, 10).Value = 49099.6  # J20
$ws.Cells.Item(20, 12).Value = 49099.6  # L20
$ws.Cells.Item(20, 14).Value = -49571.6  # N20
$ws.Cells.Item(30, 8).Value = 49099.6  # H30
$ws.Cells.Item(30, 10).Value = 49099.6  # J30
$ws.Cells.Item(30, 12).Value = 49099.6  # L30
$ws.Cells.Item(30, 14).Value = -49281.6  # N30
$ws.Cells.Item(31, 8).Value = 13335203  # H31
$ws.Cells.Item(31, 9).Value = 1508.5883  # I31
$ws.Cells.Item(31, 10).Value = 41669304  # J31
$ws.Cells.Item(31, 11).Value = 1508.5883  # K31
$ws.Cells.Item(31, 12).Value = 41669304  # L31
$ws.Cells.Item(31, 13).Value = -1213.5883  # M31
$ws.Cells.Item(31, 14).Value = -41669894  # N31
$ws.Cells.Item(34, 8).Value = 13335203  # H34
$ws.Cells.Item(34, 9).Value = 1508.5883  # I34
$ws.Cells.Item(34, 10).Value = 41669304  # J34
$ws.Cells.Item(34, 11).Value = 1508.5883  # K34
$ws.Cells.Item(34, 12).Value = 41669304  # L34
$ws.Cells.Item(34, 13).Value = -1306.5883  # M34
$ws.Cells.Item(34, 14).Value = -41669708  # N34
$ws.Cells.Item(50, 8).Value = 9292.333000000001  # H50
$ws.Cells.Item(50, 10).Value = 9292.333000000001  # J50
$ws.Cells.Item(50, 12).Value = 9292.333000000001  # L50
$ws.Cells.Item(50, 14).Value = -10542.333  # N50
$ws.Cells.Item(99, 8).Value = 1603.3334  # H99
$ws.Cells.Item(99, 9).Value = 1608.5714  # I99
$ws.Cells.Item(99, 11).Value = 1608.5714  # K99
$ws.Cells.Item(99, 13).Value = -110.5714  # M99
$ws.Cells.Item(108, 8).Value = 23693  # H108
$ws.Cells.Item(108, 10).Value = 26616.25  # J108
$ws.Cells.Item(108, 12).Value = 26616.25  # L108
$ws.Cells.Item(108, 14).Value = -34296.25  # N108
$ws.Cells.Item(109, 8).Value = 10225  # H109
$ws.Cells.Item(109, 10).Value = 10257.143  # J109
$ws.Cells.Item(109, 12).Value = 10257.143  # L109
$ws.Cells.Item(109, 14).Value = -12337.143  # N109
$ws.Cells.Item(126, 8).Value = 1603.3334  # H126
$ws.Cells.Item(126, 9).Value = 1608.5714  # I126
$ws.Cells.Item(126, 11).Value = 4825.7142  # K126
$ws.Cells.Item(126, 13).Value = -2355.7142  # M126
$ws.Cells.Item(127, 8).Value = 54028  # H127
$ws.Cells.Item(127, 10).Value = 54028  # J127
$ws.Cells.Item(127, 12).Value = 54028  # L127
$ws.Cells.Item(127, 14).Value = -63948  # N127
$ws.Cells.Item(128, 8).Value = 49099.6  # H128
$ws.Cells.Item(128, 10).Value = 49099.6  # J128
$ws.Cells.Item(128, 12).Value = 49099.6  # L128
$ws.Cells.Item(128, 14).Value = -59059.6  # N128
$ws.Cells.Item(130, 8).Value = 50693.332  # H130
$ws.Cells.Item(130, 10).Value = 50693.332  # J130
$ws.Cells.Item(130, 12).Value = 50693.332  # L130
$ws.Cells.Item(130, 14).Value = -60733.332  # N130
$ws.Cells.Item(132, 8).Value = 2942.8572  # H132
$ws.Cells.Item(132, 9).Value = 2281.2273  # I132
$ws.Cells.Item(132, 10).Value = 5368.8335  # J132
$ws.Cells.Item(132, 11).Value = 6843.6819  # K132
$ws.Cells.Item(132, 12).Value = 16106.5005  # L132
$ws.Cells.Item(132, 13).Value = -4313.6819  # M132
$ws.Cells.Item(132, 14).Value = -21166.5005  # N132
$ws.Cells.Item(134, 8).Value = 3512654.2  # H134
$ws.Cells.Item(134, 9).Value = 4013.327  # I134
$ws.Cells.Item(134, 10).Value = 40002520  # J134
$ws.Cells.Item(134, 11).Value = 12039.981  # K134
$ws.Cells.Item(134, 12).Value = 120007560  # L134
$ws.Cells.Item(134, 13).Value = -9504.981  # M134
$ws.Cells.Item(134, 14).Value = -120012630  # N134
$ws.Cells.Item(135, 8).Value = 49281.188  # H135
$ws.Cells.Item(135, 10).Value = 49281.188  # J135
$ws.Cells.Item(135, 12).Value = 49281.188  # L135
$ws.Cells.Item(135, 14).Value = -59421.188  # N135

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 38400  # H57
$ws.Cells.Item(57, 10).Value = 38400  # J57
$ws.Cells.Item(57, 12).Value = 38400  # L57
$ws.Cells.Item(57, 14).Value = -40040  # N57
$ws.Cells.Item(124, 8).Value = 0  # H124
$ws.Cells.Item(124, 10).Value = 0  # J124
$ws.Cells.Item(124, 12).Value = 0  # L124
$ws.Cells.Item(124, 14).ClearContents()  # N124
$ws.Cells.Item(130, 8).Value = 48416  # H130
$ws.Cells.Item(130, 10).Value = 48416  # J130
$ws.Cells.Item(130, 12).Value = 48416  # L130
$ws.Cells.Item(130, 14).Value = -58456  # N130
$ws.Cells.Item(132, 8).Value = 2817.9092  # H132
$ws.Cells.Item(132, 9).Value = 2164.182  # I132
$ws.Cells.Item(132, 10).Value = 3471.6365  # J132
$ws.Cells.Item(132, 11).Value = 6492.545999999999  # K132
$ws.Cells.Item(132, 12).Value = 10414.9095  # L132
$ws.Cells.Item(132, 13).Value = -3962.545999999999  # M132
$ws.Cells.Item(132, 14).Value = -15474.9095  # N132
$ws.Cells.Item(133, 8).Value = 51139.5  # H133
$ws.Cells.Item(133, 10).Value = 51139.5  # J133
$ws.Cells.Item(133, 12).Value = 51139.5  # L133
$ws.Cells.Item(133, 14).Value = -61259.5  # N133
$ws.Cells.Item(135, 8).Value = 57764.445  # H135
$ws.Cells.Item(135, 10).Value = 57764.445  # J135
$ws.Cells.Item(135, 12).Value = 57764.445  # L135
$ws.Cells.Item(135, 14).Value = -67904.44500000001  # N135

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(92, 8).Value = 38672  # H92
$ws.Cells.Item(92, 10).Value = 38672  # J92
$ws.Cells.Item(92, 12).Value = 38672  # L92
$ws.Cells.Item(92, 14).Value = -43664  # N92
$ws.Cells.Item(127, 8).Value = 43597  # H127
$ws.Cells.Item(127, 10).Value = 43597  # J127
$ws.Cells.Item(127, 12).Value = 43597  # L127
$ws.Cells.Item(127, 14).Value = -53517  # N127
$ws.Cells.Item(130, 8).Value = 44558.168  # H130
$ws.Cells.Item(130, 10).Value = 44558.168  # J130
$ws.Cells.Item(130, 12).Value = 44558.168  # L130
$ws.Cells.Item(130, 14).Value = -54598.168  # N130
$ws.Cells.Item(132, 8).Value = 2872.28  # H132
$ws.Cells.Item(132, 9).Value = 2520.2666  # I132
$ws.Cells.Item(132, 10).Value = 3400.3  # J132
$ws.Cells.Item(132, 11).Value = 7560.7998  # K132
$ws.Cells.Item(132, 12).Value = 10200.9  # L132
$ws.Cells.Item(132, 13).Value = -5030.7998  # M132
$ws.Cells.Item(132, 14).Value = -15260.9  # N132
$ws.Cells.Item(133, 8).Value = 46268  # H133
$ws.Cells.Item(133, 10).Value = 46268  # J133
$ws.Cells.Item(133, 12).Value = 46268  # L133
$ws.Cells.Item(133, 14).Value = -51328  # N133

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(93, 8).Value = 27222.25  # H93
$ws.Cells.Item(93, 10).Value = 27222.25  # J93
$ws.Cells.Item(93, 12).Value = 27222.25  # L93
$ws.Cells.Item(93, 14).Value = -32214.25  # N93
$ws.Cells.Item(109, 8).Value = 5800  # H109
$ws.Cells.Item(109, 10).Value = 5800  # J109
$ws.Cells.Item(109, 12).Value = 5800  # L109
$ws.Cells.Item(109, 14).Value = -8574  # N109
$ws.Cells.Item(126, 8).Value = 1606.7742  # H126
$ws.Cells.Item(126, 9).Value = 1533.125  # I126
$ws.Cells.Item(126, 10).Value = 1859.2858  # J126
$ws.Cells.Item(126, 11).Value = 4599.375  # K126
$ws.Cells.Item(126, 12).Value = 5577.857400000001  # L126
$ws.Cells.Item(126, 13).Value = -2129.375  # M126
$ws.Cells.Item(126, 14).Value = -10517.8574  # N126
$ws.Cells.Item(128, 8).Value = 49995  # H128
$ws.Cells.Item(128, 10).Value = 49995  # J128
$ws.Cells.Item(128, 12).Value = 49995  # L128
$ws.Cells.Item(128, 14).Value = -59955  # N128
$ws.Cells.Item(132, 8).Value = 1882.7297  # H132
$ws.Cells.Item(132, 9).Value = 1899.7693  # I132
$ws.Cells.Item(132, 10).Value = 1842.4546  # J132
$ws.Cells.Item(132, 11).Value = 5699.3079  # K132
$ws.Cells.Item(132, 12).Value = 5527.3638  # L132
$ws.Cells.Item(132, 13).Value = -3169.3079  # M132
$ws.Cells.Item(132, 14).Value = -10587.3638  # N132

Write-Output "Applied 213 cell updates"